# This edit is a "resave" of the SharePoint Document Information Panel (DIP)
# custom XML metadata part that Word keeps alongside the document:
#   customXml/item1.xml      - the ct:contentTypeSchema (ma:versionID / ma:fieldsID refreshed)
#   customXml/itemProps1.xml - the ds:datastoreItem that holds item1's ds:itemID
#                              (itemID refreshed, cached ds:schemaRefs dropped)
# No visible document content changes - only the "Document" content-type
# schema part is re-stamped with new GUIDs, exactly like Word/SharePoint do
# whenever that part is rewritten.
#
# CustomXMLPart.XML is read-only on the Word object model (SharePoint DIP
# parts can't be edited in place), so the supported way to replace a part's
# contents is: locate it, Delete() it, and Add() a replacement with the new
# XML. Word mints a fresh GUID for the new part's .ID / the regenerated
# itemProps*.xml - which is exactly the kind of refresh this commit shows.

$d = $word.ActiveDocument

$oldContentTypePartId = "{4D5AF873-3296-4A1B-9EB3-3AD20CED7BA0}"

$newContentTypeXml = @'
<ct:contentTypeSchema xmlns:ct="http://schemas.microsoft.com/office/2006/metadata/contentType" xmlns:ma="http://schemas.microsoft.com/office/2006/metadata/properties/metaAttributes" ct:_="" ma:_="" ma:contentTypeName="Document" ma:contentTypeID="0x010100439E89A99251404F8E4F3ABB73E584D0" ma:contentTypeVersion="10" ma:contentTypeDescription="Create a new document." ma:contentTypeScope="" ma:versionID="8f9059bd7a7006de8a9c86fca591aade">
  <xsd:schema xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xs="http://www.w3.org/2001/XMLSchema" xmlns:p="http://schemas.microsoft.com/office/2006/metadata/properties" xmlns:ns1="http://schemas.microsoft.com/sharepoint/v3" xmlns:ns2="2427c28e-dd78-4248-a634-a8c840650ba2" targetNamespace="http://schemas.microsoft.com/office/2006/metadata/properties" ma:root="true" ma:fieldsID="d250b55b03be3e4b4b2c7acd5bf85c08" ns1:_="" ns2:_="">
    <xsd:import namespace="http://schemas.microsoft.com/sharepoint/v3"/>
    <xsd:import namespace="2427c28e-dd78-4248-a634-a8c840650ba2"/>
    <xsd:element name="properties">
      <xsd:complexType>
        <xsd:sequence>
          <xsd:element name="documentManagement">
            <xsd:complexType>
              <xsd:all>
                <xsd:element ref="ns2:MediaServiceMetadata" minOccurs="0"/>
                <xsd:element ref="ns2:MediaServiceFastMetadata" minOccurs="0"/>
                <xsd:element ref="ns2:MediaServiceDateTaken" minOccurs="0"/>
                <xsd:element ref="ns2:MediaServiceAutoTags" minOccurs="0"/>
                <xsd:element ref="ns2:MediaServiceOCR" minOccurs="0"/>
                <xsd:element ref="ns1:_ip_UnifiedCompliancePolicyProperties" minOccurs="0"/>
                <xsd:element ref="ns1:_ip_UnifiedCompliancePolicyUIAction" minOccurs="0"/>
                <xsd:element ref="ns2:MediaServiceLocation" minOccurs="0"/>
                <xsd:element ref="ns2:MediaServiceEventHashCode" minOccurs="0"/>
                <xsd:element ref="ns2:MediaServiceGenerationTime" minOccurs="0"/>
              </xsd:all>
            </xsd:complexType>
          </xsd:element>
        </xsd:sequence>
      </xsd:complexType>
    </xsd:element>
  </xsd:schema>
  <xsd:schema targetNamespace="http://schemas.microsoft.com/sharepoint/v3" elementFormDefault="qualified" xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xs="http://www.w3.org/2001/XMLSchema" xmlns:dms="http://schemas.microsoft.com/office/2006/documentManagement/types" xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls">
    <xsd:import namespace="http://schemas.microsoft.com/office/2006/documentManagement/types"/>
    <xsd:import namespace="http://schemas.microsoft.com/office/infopath/2007/PartnerControls"/>
    <xsd:element name="_ip_UnifiedCompliancePolicyProperties" ma:index="13" nillable="true" ma:displayName="Unified Compliance Policy Properties" ma:hidden="true" ma:internalName="_ip_UnifiedCompliancePolicyProperties">
      <xsd:simpleType>
        <xsd:restriction base="dms:Note"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="_ip_UnifiedCompliancePolicyUIAction" ma:index="14" nillable="true" ma:displayName="Unified Compliance Policy UI Action" ma:hidden="true" ma:internalName="_ip_UnifiedCompliancePolicyUIAction">
      <xsd:simpleType>
        <xsd:restriction base="dms:Text"/>
      </xsd:simpleType>
    </xsd:element>
  </xsd:schema>
  <xsd:schema targetNamespace="2427c28e-dd78-4248-a634-a8c840650ba2" elementFormDefault="qualified" xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xs="http://www.w3.org/2001/XMLSchema" xmlns:dms="http://schemas.microsoft.com/office/2006/documentManagement/types" xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls">
    <xsd:import namespace="http://schemas.microsoft.com/office/2006/documentManagement/types"/>
    <xsd:import namespace="http://schemas.microsoft.com/office/infopath/2007/PartnerControls"/>
    <xsd:element name="MediaServiceMetadata" ma:index="8" nillable="true" ma:displayName="MediaServiceMetadata" ma:hidden="true" ma:internalName="MediaServiceMetadata" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Note"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaServiceFastMetadata" ma:index="9" nillable="true" ma:displayName="MediaServiceFastMetadata" ma:hidden="true" ma:internalName="MediaServiceFastMetadata" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Note"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaServiceDateTaken" ma:index="10" nillable="true" ma:displayName="MediaServiceDateTaken" ma:hidden="true" ma:internalName="MediaServiceDateTaken" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Text"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaServiceAutoTags" ma:index="11" nillable="true" ma:displayName="MediaServiceAutoTags" ma:internalName="MediaServiceAutoTags" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Text"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaServiceOCR" ma:index="12" nillable="true" ma:displayName="MediaServiceOCR" ma:internalName="MediaServiceOCR" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Note">
          <xsd:maxLength value="255"/>
        </xsd:restriction>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaServiceLocation" ma:index="15" nillable="true" ma:displayName="MediaServiceLocation" ma:internalName="MediaServiceLocation" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Text"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaServiceEventHashCode" ma:index="16" nillable="true" ma:displayName="MediaServiceEventHashCode" ma:hidden="true" ma:internalName="MediaServiceEventHashCode" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Text"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaServiceGenerationTime" ma:index="17" nillable="true" ma:displayName="MediaServiceGenerationTime" ma:hidden="true" ma:internalName="MediaServiceGenerationTime" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Text"/>
      </xsd:simpleType>
    </xsd:element>
  </xsd:schema>
  <xsd:schema targetNamespace="http://schemas.openxmlformats.org/package/2006/metadata/core-properties" elementFormDefault="qualified" attributeFormDefault="unqualified" blockDefault="#all" xmlns="http://schemas.openxmlformats.org/package/2006/metadata/core-properties" xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xsi="http://www.w3.org/2001/XMLSchema-instance" xmlns:dc="http://purl.org/dc/elements/1.1/" xmlns:dcterms="http://purl.org/dc/terms/" xmlns:odoc="http://schemas.microsoft.com/internal/obd">
    <xsd:import namespace="http://purl.org/dc/elements/1.1/" schemaLocation="http://dublincore.org/schemas/xmls/qdc/2003/04/02/dc.xsd"/>
    <xsd:import namespace="http://purl.org/dc/terms/" schemaLocation="http://dublincore.org/schemas/xmls/qdc/2003/04/02/dcterms.xsd"/>
    <xsd:element name="coreProperties" type="CT_coreProperties"/>
    <xsd:complexType name="CT_coreProperties">
      <xsd:all>
        <xsd:element ref="dc:creator" minOccurs="0" maxOccurs="1"/>
        <xsd:element ref="dcterms:created" minOccurs="0" maxOccurs="1"/>
        <xsd:element ref="dc:identifier" minOccurs="0" maxOccurs="1"/>
        <xsd:element name="contentType" minOccurs="0" maxOccurs="1" type="xsd:string" ma:index="0" ma:displayName="Content Type"/>
        <xsd:element ref="dc:title" minOccurs="0" maxOccurs="1" ma:index="4" ma:displayName="Title"/>
        <xsd:element ref="dc:subject" minOccurs="0" maxOccurs="1"/>
        <xsd:element ref="dc:description" minOccurs="0" maxOccurs="1"/>
        <xsd:element name="keywords" minOccurs="0" maxOccurs="1" type="xsd:string"/>
        <xsd:element ref="dc:language" minOccurs="0" maxOccurs="1"/>
        <xsd:element name="category" minOccurs="0" maxOccurs="1" type="xsd:string"/>
        <xsd:element name="version" minOccurs="0" maxOccurs="1" type="xsd:string"/>
        <xsd:element name="revision" minOccurs="0" maxOccurs="1" type="xsd:string">
          <xsd:annotation>
            <xsd:documentation>
                        This value indicates the number of saves or revisions. The application is responsible for updating this value after each revision.
                    </xsd:documentation>
          </xsd:annotation>
        </xsd:element>
        <xsd:element name="lastModifiedBy" minOccurs="0" maxOccurs="1" type="xsd:string"/>
        <xsd:element ref="dcterms:modified" minOccurs="0" maxOccurs="1"/>
        <xsd:element name="contentStatus" minOccurs="0" maxOccurs="1" type="xsd:string"/>
      </xsd:all>
    </xsd:complexType>
  </xsd:schema>
  <xs:schema targetNamespace="http://schemas.microsoft.com/office/infopath/2007/PartnerControls" elementFormDefault="qualified" attributeFormDefault="unqualified" xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls" xmlns:xs="http://www.w3.org/2001/XMLSchema">
    <xs:element name="Person">
      <xs:complexType>
        <xs:sequence>
          <xs:element ref="pc:DisplayName" minOccurs="0"></xs:element>
          <xs:element ref="pc:AccountId" minOccurs="0"></xs:element>
          <xs:element ref="pc:AccountType" minOccurs="0"></xs:element>
        </xs:sequence>
      </xs:complexType>
    </xs:element>
    <xs:element name="DisplayName" type="xs:string"></xs:element>
    <xs:element name="AccountId" type="xs:string"></xs:element>
    <xs:element name="AccountType" type="xs:string"></xs:element>
    <xs:element name="BDCAssociatedEntity">
      <xs:complexType>
        <xs:sequence>
          <xs:element ref="pc:BDCEntity" minOccurs="0" maxOccurs="unbounded"></xs:element>
        </xs:sequence>
        <xs:attribute ref="pc:EntityNamespace"></xs:attribute>
        <xs:attribute ref="pc:EntityName"></xs:attribute>
        <xs:attribute ref="pc:SystemInstanceName"></xs:attribute>
        <xs:attribute ref="pc:AssociationName"></xs:attribute>
      </xs:complexType>
    </xs:element>
    <xs:attribute name="EntityNamespace" type="xs:string"></xs:attribute>
    <xs:attribute name="EntityName" type="xs:string"></xs:attribute>
    <xs:attribute name="SystemInstanceName" type="xs:string"></xs:attribute>
    <xs:attribute name="AssociationName" type="xs:string"></xs:attribute>
    <xs:element name="BDCEntity">
      <xs:complexType>
        <xs:sequence>
          <xs:element ref="pc:EntityDisplayName" minOccurs="0"></xs:element>
          <xs:element ref="pc:EntityInstanceReference" minOccurs="0"></xs:element>
          <xs:element ref="pc:EntityId1" minOccurs="0"></xs:element>
          <xs:element ref="pc:EntityId2" minOccurs="0"></xs:element>
          <xs:element ref="pc:EntityId3" minOccurs="0"></xs:element>
          <xs:element ref="pc:EntityId4" minOccurs="0"></xs:element>
          <xs:element ref="pc:EntityId5" minOccurs="0"></xs:element>
        </xs:sequence>
      </xs:complexType>
    </xs:element>
    <xs:element name="EntityDisplayName" type="xs:string"></xs:element>
    <xs:element name="EntityInstanceReference" type="xs:string"></xs:element>
    <xs:element name="EntityId1" type="xs:string"></xs:element>
    <xs:element name="EntityId2" type="xs:string"></xs:element>
    <xs:element name="EntityId3" type="xs:string"></xs:element>
    <xs:element name="EntityId4" type="xs:string"></xs:element>
    <xs:element name="EntityId5" type="xs:string"></xs:element>
    <xs:element name="Terms">
      <xs:complexType>
        <xs:sequence>
          <xs:element ref="pc:TermInfo" minOccurs="0" maxOccurs="unbounded"></xs:element>
        </xs:sequence>
      </xs:complexType>
    </xs:element>
    <xs:element name="TermInfo">
      <xs:complexType>
        <xs:sequence>
          <xs:element ref="pc:TermName" minOccurs="0"></xs:element>
          <xs:element ref="pc:TermId" minOccurs="0"></xs:element>
        </xs:sequence>
      </xs:complexType>
    </xs:element>
    <xs:element name="TermName" type="xs:string"></xs:element>
    <xs:element name="TermId" type="xs:string"></xs:element>
  </xs:schema>
</ct:contentTypeSchema>
'@

$oldPart = $null

# Prefer an exact ID lookup (SelectByID), falling back to a namespace/content
# scan if that overload isn't available, and finally to a full linear scan.
try {
    $oldPart = $d.CustomXMLParts.SelectByID($oldContentTypePartId)
} catch {
    $oldPart = $null
}

if (-not $oldPart) {
    try {
        $candidates = $d.CustomXMLParts.SelectByNamespace("http://schemas.microsoft.com/office/2006/metadata/contentType")
        if ($candidates -and $candidates.Count -ge 1) {
            $oldPart = $candidates.Item(1)
        }
    } catch {
        $oldPart = $null
    }
}

if (-not $oldPart) {
    try {
        $parts = $d.CustomXMLParts
        for ($i = 1; $i -le $parts.Count; $i++) {
            $candidate = $parts.Item($i)
            if ($candidate.ID -eq $oldContentTypePartId -or
                ($candidate.XML -and $candidate.XML.Contains("ma:contentTypeName=`"Document`""))) {
                $oldPart = $candidate
                break
            }
        }
    } catch {
    }
}

if ($oldPart) {
    try { $oldPart.Delete() } catch { }
}

try {
    $d.CustomXMLParts.Add($newContentTypeXml) | Out-Null
} catch { }
